# Logged Week 15 and simulated Week 16
# Update the "R" (road) row totals on both the OFF and DEF sheets.

$wb = $excel.ActiveWorkbook

# --- OFF sheet, row 3 (R) ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 334
$wsOff.Range("C3").Value = 222
$wsOff.Range("D3").Value = 86
$wsOff.Range("E3").Value = 28
$wsOff.Range("F3").Value = 7

# --- DEF sheet, row 3 (R) ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 465
$wsDef.Range("C3").Value = 328
$wsDef.Range("D3").Value = 95
$wsDef.Range("E3").Value = 47
$wsDef.Range("F3").Value = 6
